$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of class-specific accuracy data
$ws.Range("A4").Value = 1000
$ws.Range("B4").Value = 50
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = 30
$ws.Range("E4").Value = 2000
$ws.Range("F4").Value = 0.993
$ws.Range("G4").Value = 0.8519
$ws.Range("H4").Value = 0.9118

# Update selection to reflect the new active cell after data entry
$ws.Range("A5").Select() | Out-Null
